$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13/14 swap: Avalanche <-> Polygon (plain text fields) ---
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"

# --- E column (Volume) updates: plain text (leading/trailing spaces keep Excel
#     from re-interpreting these as numeric percentages) ---
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  -4.28%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("E11").Value = "  +6.71%  "
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("E13").Value = "  -3.94%  "
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("E17").Value = "  -2.10%  "
$ws.Range("E18").Value = "  -1.13%  "
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("E21").Value = "  +3.07%  "
$ws.Range("E22").Value = "  -3.20%  "
$ws.Range("E23").Value = "  -2.32%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("E27").Value = "  -3.93%  "
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("E30").Value = "  -7.94%  "
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("E33").Value = "  -3.17%  "
$ws.Range("E34").Value = "  +2.76%  "
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  -3.67%  "
$ws.Range("E39").Value = "  -5.93%  "
$ws.Range("E40").Value = "  -2.02%  "
$ws.Range("E41").Value = "  -0.94%  "
$ws.Range("E42").Value = "  +4.69%  "
$ws.Range("E43").Value = "  -2.44%  "
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("E45").Value = "  +0.88%  "
$ws.Range("E46").Value = "  +1.57%  "
$ws.Range("E47").Value = "  -4.75%  "
$ws.Range("E48").Value = "  -4.59%  "
$ws.Range("E49").Value = "  -4.64%  "
$ws.Range("E50").Value = "  +3.62%  "
$ws.Range("E51").Value = "  +3.94%  "

# --- D column (Price) updates -----------------------------------------------
# Many of the new prices are plain numeric-looking strings (e.g. "243.46") that
# Excel would otherwise auto-convert to a true number, losing the original
# text formatting used throughout this sheet. To keep them as text we flip the
# cell to the "@" (Text) number format immediately before assigning the value,
# then clear the (now unneeded) formatting again right away so the cell is left
# without any explicit style, just like its neighbours.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "36.407.94"
$cell.ClearFormats()
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.941.28"
$cell.ClearFormats()
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "243.46"
$cell.ClearFormats()
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "57.74"
$cell.ClearFormats()
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.367"
$cell.ClearFormats()
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "55.79"
$cell.ClearFormats()
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0836"
$cell.ClearFormats()
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.823"
$cell.ClearFormats()
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "21.54"
$cell.ClearFormats()
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "2.229.14"
$cell.ClearFormats()
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "13.59"
$cell.ClearFormats()
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "5.25"
$cell.ClearFormats()
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "1.954.13"
$cell.ClearFormats()
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "36.339.75"
$cell.ClearFormats()
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0867"
$cell.ClearFormats()
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "229.51"
$cell.ClearFormats()
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.03"
$cell.ClearFormats()
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.44"
$cell.ClearFormats()
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "9.26"
$cell.ClearFormats()
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "162.13"
$cell.ClearFormats()
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.123"
$cell.ClearFormats()
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "4.67"
$cell.ClearFormats()
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.0628"
$cell.ClearFormats()
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "6.24"
$cell.ClearFormats()
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "3.03"
$cell.ClearFormats()
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.0971"
$cell.ClearFormats()
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "2.94"
$cell.ClearFormats()
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "1.17"
$cell.ClearFormats()
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.0208"
$cell.ClearFormats()
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "16.09"
$cell.ClearFormats()
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "1.353.11"
$cell.ClearFormats()
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "87.63"
$cell.ClearFormats()
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "7.15"
$cell.ClearFormats()
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "2.81"
$cell.ClearFormats()
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "45.13"
$cell.ClearFormats()

Write-Host "Updated cryptos list."
